# dfPrediccion.xlsx update:
#  - New prediction run: fecha (column B) moves from 2024-04-02 (45384) to
#    2024-04-05 (45387) for every remaining row.
#  - prediccion (column D) values are refreshed with the new model output.
#  - Horizons 27, 28 and 29 (rows 29-31) are dropped — the sheet now only
#    covers horizon 0 through 26 (rows 2-28).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New prediction values for horizon 0..26 (rows 2..28)
$newValues = @(
    7552.7978515625,
    7415.4716796875,
    7461.67138671875,
    7363.22607421875,
    7325.51611328125,
    7363.67919921875,
    7235.015625,
    7426.04296875,
    7442.71142578125,
    7461.1865234375,
    7463.392578125,
    7494.93896484375,
    7428.25390625,
    7493.05517578125,
    7434.95166015625,
    7498.55126953125,
    7515.53369140625,
    7424.79150390625,
    7451.2919921875,
    7428.0908203125,
    7527.73828125,
    7499.58203125,
    7473.0791015625,
    7421.859375,
    7439.10302734375,
    7440.11572265625,
    7376.45751953125
)

$newDate = 45387

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $newDate
    $ws.Cells.Item($row, 4).Value = $newValues[$i]
}

# Drop the trailing horizons (27, 28, 29) that no longer exist — rows 29-31.
$ws.Range("A29:D31").Delete()

Write-Output "Updated dfPrediccion sheet: new fecha=$newDate, rows now 1:28"
